$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (price + 1h volume-change columns) with the latest
# scrape. A couple of rows (Bittensor / Binance-PegBSC-USD) also swapped rank
# order, so their Coin name / Link / Price / Volume cells are rewritten too.
# Price values that look like plain numbers ("553.34", "1.00", ...) are
# prefixed with a leading apostrophe so Excel keeps storing them as text
# (matching the original inline-string cell type) instead of silently
# re-interpreting them as numeric values.

$ws.Range("D2").Value = '61.478.28'
$ws.Range("E2").Value = '  -3.55%  '
$ws.Range("D3").Value = '2.484.43'
$ws.Range("E3").Value = '  -5.84%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''553.34'
$ws.Range("E5").Value = '  -4.51%  '
$ws.Range("D6").Value = '''146.45'
$ws.Range("E6").Value = '  -5.70%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.601'
$ws.Range("E8").Value = '  -3.06%  '
$ws.Range("D9").Value = '2.486.13'
$ws.Range("E9").Value = '  -5.70%  '
$ws.Range("E10").Value = '  -8.70%  '
$ws.Range("D11").Value = '''5.45'
$ws.Range("E11").Value = '  -6.25%  '
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("E13").Value = '  -6.34%  '
$ws.Range("D14").Value = '''26.20'
$ws.Range("E14").Value = '  -7.80%  '
$ws.Range("D15").Value = '2.933.34'
$ws.Range("E15").Value = '  -5.85%  '
$ws.Range("E16").Value = '  -8.43%  '
$ws.Range("D17").Value = '61.456.63'
$ws.Range("E17").Value = '  -3.53%  '
$ws.Range("D18").Value = '2.489.50'
$ws.Range("E18").Value = '  -5.99%  '
$ws.Range("E19").Value = '  -7.70%  '
$ws.Range("D20").Value = '''7.06'
$ws.Range("E20").Value = '  -7.62%  '
$ws.Range("D21").Value = '''4.20'
$ws.Range("E21").Value = '  -6.99%  '
$ws.Range("D22").Value = '''322.40'
$ws.Range("E22").Value = '  -6.31%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '''1.88'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").Value = '''64.14'
$ws.Range("E25").Value = '  -5.71%  '
$ws.Range("D26").Value = '''0.0000100'
$ws.Range("E26").Value = '  -8.05%  '
$ws.Range("D27").Value = '2.610.08'
$ws.Range("E27").Value = '  -5.72%  '
$ws.Range("D28").Value = '''1.51'
$ws.Range("E28").Value = '  -5.54%  '
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").Value = '''544.68'
$ws.Range("E29").Value = '  -9.57%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  -9.52%  '
$ws.Range("D32").Value = '''7.78'
$ws.Range("E32").Value = '  -4.68%  '
$ws.Range("E33").Value = '  -5.75%  '
$ws.Range("E34").Value = '  -7.33%  '
$ws.Range("D35").Value = '''1.60'
$ws.Range("E35").Value = '  -8.06%  '
$ws.Range("D36").Value = '''5.90'
$ws.Range("E36").Value = '  -10.07%  '
$ws.Range("D37").Value = '''4.88'
$ws.Range("E37").Value = '  -10.18%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("E39").Value = '  -5.16%  '
$ws.Range("D40").Value = '''18.59'
$ws.Range("E40").Value = '  -5.72%  '
$ws.Range("D41").Value = '''148.39'
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("E42").Value = '  -8.37%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '''40.42'
$ws.Range("E44").Value = '  -3.51%  '
$ws.Range("D45").Value = '''2.38'
$ws.Range("E45").Value = '  -6.79%  '
$ws.Range("D46").Value = '''147.80'
$ws.Range("E46").Value = '  -8.00%  '
$ws.Range("E47").Value = '  -6.70%  '
$ws.Range("D48").Value = '''21.19'
$ws.Range("E48").Value = '  -12.36%  '
$ws.Range("D49").Value = '''0.0541'
$ws.Range("E49").Value = '  -7.46%  '
$ws.Range("D50").Value = '''0.600'
$ws.Range("E50").Value = '  -5.16%  '
$ws.Range("D51").Value = '''0.0949'
$ws.Range("E51").Value = '  -4.69%  '
